$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6600
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496

$ws.Range("H67").Value = 6600
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716

$ws.Range("H74").Value = 9999
$ws.Range("I74").Value = 9999
$ws.Range("J74").Value = 9999
$ws.Range("K74").Value = 9999
$ws.Range("L74").Value = 9999
$ws.Range("M74").Value = -9063
$ws.Range("N74").Value = -11871

$ws.Range("H77").Value = 9999
$ws.Range("I77").Value = 9999
$ws.Range("J77").Value = 9999
$ws.Range("K77").Value = 49995
$ws.Range("L77").Value = 49995
$ws.Range("M77").Value = -45315
$ws.Range("N77").Value = -59355

$ws.Range("H95").Value = 38000
$ws.Range("J95").Value = 38000
$ws.Range("L95").Value = 38000
$ws.Range("N95").Value = -43492

$ws.Range("H132").Value = 2003.8667
$ws.Range("J132").Value = 400
$ws.Range("L132").Value = 1200
$ws.Range("N132").Value = -6260

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 666
$ws.Range("I3").Value = 666
$ws.Range("K3").Value = 666
$ws.Range("M3").Value = -551

$ws.Range("H5").Value = 300.33334
$ws.Range("I5").Value = 400.5
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 400.5
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -288.5
$ws.Range("N5").Value = -324

$ws.Range("H22").Value = 3083
$ws.Range("I22").Value = 2110.6667
$ws.Range("K22").Value = 2110.6667
$ws.Range("M22").Value = -1811.6667

$ws.Range("H32").Value = 8360
$ws.Range("I32").Value = 6551.636
$ws.Range("K32").Value = 6551.636
$ws.Range("M32").Value = -6264.636

$ws.Range("H41").Value = 967.2857
$ws.Range("I41").Value = 651.5
$ws.Range("J41").Value = 2862
$ws.Range("K41").Value = 651.5
$ws.Range("L41").Value = 2862
$ws.Range("M41").Value = -237.5
$ws.Range("N41").Value = -3690

$ws.Range("H88").Value = 2583
$ws.Range("I88").Value = 2449.5
$ws.Range("J88").Value = 2649.75
$ws.Range("K88").Value = 2449.5
$ws.Range("L88").Value = 2649.75
$ws.Range("M88").Value = -2043.5
$ws.Range("N88").Value = -3461.75

$ws.Range("H91").Value = 2583
$ws.Range("I91").Value = 2449.5
$ws.Range("J91").Value = 2649.75
$ws.Range("K91").Value = 2449.5
$ws.Range("L91").Value = 2649.75
$ws.Range("M91").Value = -1045.5
$ws.Range("N91").Value = -5457.75

$ws.Range("H92").Value = 137119.75
$ws.Range("J92").Value = 137119.75
$ws.Range("L92").Value = 137119.75
$ws.Range("N92").Value = -142111.75

$ws.Range("H96").Value = 19106.428
$ws.Range("J96").Value = 19106.428
$ws.Range("L96").Value = 19106.428
$ws.Range("N96").Value = -24598.428

$ws.Range("H122").Value = 5136.3076
$ws.Range("I122").Value = 2096.5
$ws.Range("K122").Value = 6289.5
$ws.Range("M122").Value = -3839.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 300.33334
$ws.Range("I4").Value = 400.5
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 400.5
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -285.5
$ws.Range("N4").Value = -330

$ws.Range("H82").Value = 46000
$ws.Range("I82").Value = 28400
$ws.Range("J82").Value = 90000
$ws.Range("K82").Value = 28400
$ws.Range("L82").Value = 90000
$ws.Range("M82").Value = -28017
$ws.Range("N82").Value = -90766

$ws.Range("H85").Value = 46000
$ws.Range("I85").Value = 28400
$ws.Range("J85").Value = 90000
$ws.Range("K85").Value = 28400
$ws.Range("L85").Value = 90000
$ws.Range("M85").Value = -27074
$ws.Range("N85").Value = -92652

$ws.Range("H86").Value = 2065.087
$ws.Range("I86").Value = 1662.6842
$ws.Range("K86").Value = 1662.6842
$ws.Range("M86").Value = -539.6841999999999

$ws.Range("H89").Value = 2065.087
$ws.Range("I89").Value = 1662.6842
$ws.Range("K89").Value = 8313.421
$ws.Range("M89").Value = -2697.421

$ws.Range("H105").Value = 1950
$ws.Range("I105").Value = 1950
$ws.Range("K105").Value = 1950
$ws.Range("M105").Value = -203

$ws.Range("H134").Value = 3033.3333
$ws.Range("I134").Value = 3033.3333
$ws.Range("K134").Value = 9099.999899999999
$ws.Range("M134").Value = -6564.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 49200
$ws.Range("J74").Value = 80000
$ws.Range("L74").Value = 80000
$ws.Range("N74").Value = -81748

$ws.Range("H77").Value = 49200
$ws.Range("J77").Value = 80000
$ws.Range("L77").Value = 240000
$ws.Range("N77").Value = -248736

$ws.Range("H141").Value = 522962.34
$ws.Range("J141").Value = 522962.34
$ws.Range("L141").Value = 522962.34
$ws.Range("N141").Value = -533322.3400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 235.66667
$ws.Range("I7").Value = 210.33333
$ws.Range("J7").Value = 286.33334
$ws.Range("K7").Value = 630.99999
$ws.Range("L7").Value = 859.0000200000001
$ws.Range("M7").Value = -518.99999
$ws.Range("N7").Value = -1083.00002

$ws.Range("H34").Value = 3223.4546
$ws.Range("J34").Value = 3622.5557
$ws.Range("L34").Value = 10867.6671
$ws.Range("N34").Value = -11035.6671

$ws.Range("H108").Value = 410.6
$ws.Range("I108").Value = 410.6
$ws.Range("K108").Value = 1231.8
$ws.Range("M108").Value = 1648.2

$ws.Range("H131").Value = 947.125
$ws.Range("I131").Value = 947.125
$ws.Range("K131").Value = 2841.375
$ws.Range("M131").Value = 2198.625

$ws.Range("H139").Value = 4027.889
$ws.Range("I139").Value = 3968.875
$ws.Range("K139").Value = 11906.625
$ws.Range("M139").Value = -6766.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 215.6
$ws.Range("I2").Value = 391.5
$ws.Range("J2").Value = 98.333336
$ws.Range("K2").Value = 391.5
$ws.Range("L2").Value = 98.333336
$ws.Range("M2").Value = -278.5
$ws.Range("N2").Value = -324.333336

$ws.Range("H3").Value = 4300
$ws.Range("I3").Value = 4300
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4300
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4184
$ws.Range("N3").ClearContents()

$ws.Range("H70").Value = 166670340
$ws.Range("I70").Value = 5504
$ws.Range("K70").Value = 5504
$ws.Range("M70").Value = -5234

$ws.Range("H73").Value = 166670340
$ws.Range("I73").Value = 5504
$ws.Range("K73").Value = 5504
$ws.Range("M73").Value = -4568

$ws.Range("H92").Value = 13889.8
$ws.Range("J92").Value = 17112.25
$ws.Range("L92").Value = 17112.25
$ws.Range("N92").Value = -20856.25

$ws.Range("H93").Value = 10000
$ws.Range("J93").Value = 10000
$ws.Range("L93").Value = 10000
$ws.Range("N93").Value = -13744

$ws.Range("H95").Value = 29625
$ws.Range("J95").Value = 29625
$ws.Range("L95").Value = 29625
$ws.Range("N95").Value = -35117

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 750.5
$ws.Range("I16").Value = 750.5
$ws.Range("K16").Value = 750.5
$ws.Range("M16").Value = -580.5

$ws.Range("H61").Value = 2217.3333
$ws.Range("I61").Value = 2076.25
$ws.Range("J61").Value = 2499.5
$ws.Range("K61").Value = 2076.25
$ws.Range("L61").Value = 2499.5
$ws.Range("M61").Value = -1874.25
$ws.Range("N61").Value = -2903.5

$ws.Range("H103").Value = 21998
$ws.Range("J103").Value = 21998
$ws.Range("L103").Value = 21998
$ws.Range("N103").Value = -24342

$ws.Range("H106").Value = 14323.75
$ws.Range("J106").Value = 14323.75
$ws.Range("L106").Value = 14323.75
$ws.Range("N106").Value = -16847.75

$ws.Range("H113").Value = 2217.3333
$ws.Range("I113").Value = 2076.25
$ws.Range("J113").Value = 2499.5
$ws.Range("K113").Value = 2076.25
$ws.Range("L113").Value = 2499.5
$ws.Range("M113").Value = 93.75
$ws.Range("N113").Value = -6839.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 21704.5
$ws.Range("J117").Value = 21704.5
$ws.Range("L117").Value = 21704.5
$ws.Range("N117").Value = -30882.5

$ws.Range("H136").Value = 2259.0715
$ws.Range("I136").Value = 2202.077
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6606.231000000001
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -4056.231000000001
$ws.Range("N136").Value = -14100
